# Apply scheduled-runner price/profit refresh to the Sheets workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8085545
$ws.Range("I116").Value = 4169242
$ws.Range("J116").Value = 12905610
$ws.Range("K116").Value = 4169242
$ws.Range("L116").Value = 12905610
$ws.Range("M116").Value = -4165800
$ws.Range("N116").Value = -12912494

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11056.029
$ws.Range("I2").Value = 14236.654
$ws.Range("J2").Value = 719
$ws.Range("K2").Value = 14236.654
$ws.Range("L2").Value = 719
$ws.Range("M2").Value = -14123.654
$ws.Range("N2").Value = -945
$ws.Range("H37").Value = 8062.7334
$ws.Range("J37").Value = 10863
$ws.Range("L37").Value = 10863
$ws.Range("N37").Value = -11409
$ws.Range("H45").Value = 2398.7273
$ws.Range("I45").Value = 1357.2727
$ws.Range("K45").Value = 1357.2727
$ws.Range("M45").Value = -980.2727
$ws.Range("H88").Value = 10393.714
$ws.Range("I88").Value = 10975
$ws.Range("J88").Value = 10161.2
$ws.Range("K88").Value = 10975
$ws.Range("L88").Value = 10161.2
$ws.Range("M88").Value = -10569
$ws.Range("N88").Value = -10973.2
$ws.Range("H91").Value = 10393.714
$ws.Range("I91").Value = 10975
$ws.Range("J91").Value = 10161.2
$ws.Range("K91").Value = 10975
$ws.Range("L91").Value = 10161.2
$ws.Range("M91").Value = -9571
$ws.Range("N91").Value = -12969.2
$ws.Range("H102").Value = 2378.5715
$ws.Range("I102").Value = 1425
$ws.Range("J102").Value = 2760
$ws.Range("K102").Value = 1425
$ws.Range("L102").Value = 2760
$ws.Range("M102").Value = 197
$ws.Range("N102").Value = -6004
$ws.Range("H116").Value = 11056.029
$ws.Range("I116").Value = 14236.654
$ws.Range("J116").Value = 719
$ws.Range("K116").Value = 14236.654
$ws.Range("L116").Value = 719
$ws.Range("M116").Value = -11942.654
$ws.Range("N116").Value = -5307
$ws.Range("H132").Value = 15615619
$ws.Range("I132").Value = 22708376
$ws.Range("J132").Value = 2316700
$ws.Range("K132").Value = 68125128
$ws.Range("L132").Value = 6950100
$ws.Range("M132").Value = -68122598
$ws.Range("N132").Value = -6955160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11056.029
$ws.Range("I3").Value = 14236.654
$ws.Range("J3").Value = 719
$ws.Range("K3").Value = 14236.654
$ws.Range("L3").Value = 719
$ws.Range("M3").Value = -14122.654
$ws.Range("N3").Value = -947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1392421.5
$ws.Range("I31").Value = 3206244
$ws.Range("K31").Value = 3206244
$ws.Range("M31").Value = -3205949
$ws.Range("H34").Value = 1392421.5
$ws.Range("I34").Value = 3206244
$ws.Range("K34").Value = 3206244
$ws.Range("M34").Value = -3206042
$ws.Range("H50").Value = 12423.5
$ws.Range("J50").Value = 12423.5
$ws.Range("L50").Value = 12423.5
$ws.Range("N50").Value = -13673.5
$ws.Range("H51").Value = 9260.6
$ws.Range("J51").Value = 9260.6
$ws.Range("L51").Value = 9260.6
$ws.Range("N51").Value = -10732.6
$ws.Range("H58").Value = 2532732.5
$ws.Range("I58").Value = 12225.223
$ws.Range("J58").Value = 5053240
$ws.Range("K58").Value = 12225.223
$ws.Range("L58").Value = 5053240
$ws.Range("M58").Value = -12022.223
$ws.Range("N58").Value = -5053646
$ws.Range("H59").Value = 16500.5
$ws.Range("J59").Value = 16500.5
$ws.Range("L59").Value = 16500.5
$ws.Range("N59").Value = -18790.5
$ws.Range("H60").Value = 4500.3
$ws.Range("I60").Value = 1966.6666
$ws.Range("K60").Value = 1966.6666
$ws.Range("M60").Value = -1455.6666
$ws.Range("H61").Value = 9260.6
$ws.Range("J61").Value = 9260.6
$ws.Range("L61").Value = 9260.6
$ws.Range("N61").Value = -9956.6
$ws.Range("H68").Value = 17849.8
$ws.Range("J68").Value = 17849.8
$ws.Range("L68").Value = 17849.8
$ws.Range("N68").Value = -19347.8
$ws.Range("H71").Value = 17849.8
$ws.Range("J71").Value = 17849.8
$ws.Range("L71").Value = 53549.39999999999
$ws.Range("N71").Value = -61037.39999999999
$ws.Range("H74").Value = 19391.445
$ws.Range("J74").Value = 19391.445
$ws.Range("L74").Value = 19391.445
$ws.Range("N74").Value = -21139.445
$ws.Range("H77").Value = 19391.445
$ws.Range("J77").Value = 19391.445
$ws.Range("L77").Value = 58174.335
$ws.Range("N77").Value = -66910.33499999999
$ws.Range("I94").Value = 798
$ws.Range("J94").Value = 45463100
$ws.Range("K94").Value = 798
$ws.Range("L94").Value = 45463100
$ws.Range("M94").Value = -347
$ws.Range("N94").Value = -45464002
$ws.Range("H99").Value = 16608.143
$ws.Range("I99").Value = 12428.571
$ws.Range("J99").Value = 20787.715
$ws.Range("K99").Value = 12428.571
$ws.Range("L99").Value = 20787.715
$ws.Range("M99").Value = -10930.571
$ws.Range("N99").Value = -23783.715
$ws.Range("H126").Value = 16608.143
$ws.Range("I126").Value = 12428.571
$ws.Range("J126").Value = 20787.715
$ws.Range("K126").Value = 37285.713
$ws.Range("L126").Value = 62363.145
$ws.Range("M126").Value = -34815.713
$ws.Range("N126").Value = -67303.145
$ws.Range("H132").Value = 1942.5
$ws.Range("I132").Value = 1101.5
$ws.Range("K132").Value = 3304.5
$ws.Range("M132").Value = -774.5
$ws.Range("H134").Value = 1820135.4
$ws.Range("I134").Value = 2148.3572
$ws.Range("J134").Value = 5001612.5
$ws.Range("K134").Value = 6445.071599999999
$ws.Range("L134").Value = 15004837.5
$ws.Range("M134").Value = -3910.071599999999
$ws.Range("N134").Value = -15009907.5
$ws.Range("H136").Value = 2532732.5
$ws.Range("I136").Value = 12225.223
$ws.Range("J136").Value = 5053240
$ws.Range("K136").Value = 36675.669
$ws.Range("L136").Value = 15159720
$ws.Range("M136").Value = -34125.669
$ws.Range("N136").Value = -15164820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 868.4299999999999
$ws.Range("I131").Value = 295
$ws.Range("J131").Value = 905.0319
$ws.Range("K131").Value = 885
$ws.Range("L131").Value = 2715.0957
$ws.Range("M131").Value = 4155
$ws.Range("N131").Value = -12795.0957

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 8853.536
$ws.Range("I80").Value = 4614.2856
$ws.Range("J80").Value = 13092.786
$ws.Range("K80").Value = 4614.2856
$ws.Range("L80").Value = 13092.786
$ws.Range("M80").Value = -3616.2856
$ws.Range("N80").Value = -15088.786
$ws.Range("H83").Value = 8853.536
$ws.Range("I83").Value = 4614.2856
$ws.Range("J83").Value = 13092.786
$ws.Range("K83").Value = 23071.428
$ws.Range("L83").Value = 65463.93
$ws.Range("M83").Value = -18079.428
$ws.Range("N83").Value = -75447.92999999999
$ws.Range("H132").Value = 8465226
$ws.Range("I132").Value = 11256834
$ws.Range("J132").Value = 5053260
$ws.Range("K132").Value = 33770502
$ws.Range("L132").Value = 15159780
$ws.Range("M132").Value = -33767972
$ws.Range("N132").Value = -15164840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5814670
$ws.Range("I46").Value = 723.0133
$ws.Range("J46").Value = 45455216
$ws.Range("K46").Value = 723.0133
$ws.Range("L46").Value = 45455216
$ws.Range("M46").Value = -535.0133
$ws.Range("N46").Value = -45455592
$ws.Range("H68").Value = 2911.8462
$ws.Range("I68").Value = 1345.1428
$ws.Range("J68").Value = 4739.6665
$ws.Range("K68").Value = 1345.1428
$ws.Range("L68").Value = 4739.6665
$ws.Range("M68").Value = -596.1428000000001
$ws.Range("N68").Value = -6237.6665
$ws.Range("H71").Value = 2911.8462
$ws.Range("I71").Value = 1345.1428
$ws.Range("J71").Value = 4739.6665
$ws.Range("K71").Value = 6725.714
$ws.Range("L71").Value = 23698.3325
$ws.Range("M71").Value = -2981.714
$ws.Range("N71").Value = -31186.3325
$ws.Range("H82").Value = 3406.24
$ws.Range("I82").Value = 1038.1
$ws.Range("J82").Value = 4985
$ws.Range("K82").Value = 1038.1
$ws.Range("L82").Value = 4985
$ws.Range("M82").Value = -677.0999999999999
$ws.Range("N82").Value = -5707
$ws.Range("H85").Value = 3406.24
$ws.Range("I85").Value = 1038.1
$ws.Range("J85").Value = 4985
$ws.Range("K85").Value = 1038.1
$ws.Range("L85").Value = 4985
$ws.Range("M85").Value = 209.9000000000001
$ws.Range("N85").Value = -7481
